# Sprint backlog cleanup: the "D1-Sprint 2" tab was an accidental duplicate
# of "D1-Sprint 1" (identical contents). Remove it, then renumber the
# remaining D2/D3 sprint tabs so the sprint numbers are contiguous again.

$wb = $excel.ActiveWorkbook

# Excel requires at least one sheet to remain visible while deleting, and
# normally prompts to confirm deleting a sheet; DisplayAlerts off mirrors
# clicking "Delete" in that dialog.
$excel.DisplayAlerts = $false

# Remove the duplicate "D1-Sprint 2" worksheet entirely.
$wb.Worksheets.Item("D1-Sprint 2").Delete() | Out-Null

$excel.DisplayAlerts = $true

# Renumber / rename the remaining sprint tabs so the sequence is contiguous.
$wsD2S3 = $wb.Worksheets.Item("D2-Sprint 3")
$wsD2S3.Range("D2").Value = 2
$wsD2S3.Name = "D2-Sprint 2"

$wsD2S4 = $wb.Worksheets.Item("D2-Sprint 4")
$wsD2S4.Range("D2").Value = 3
$wsD2S4.Name = "D2-Sprint 3"

$wsD3S5 = $wb.Worksheets.Item("D3-Sprint 5")
$wsD3S5.Name = "D3-Sprint 4"

$wsD3S6 = $wb.Worksheets.Item("D3-Sprint 6")
$wsD3S6.Name = "D3-Sprint 5"
$wsD3S6.Range("H19").Select() | Out-Null
